$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Ativo"
[void]$ws.Range("B4").Select()
